# Se modifica y valida la implementacion del caso de uso registrar pago alumno
# Se corrige la seleccion de promociones en el caso de uso registrar pago de alumno

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# --- Fill in the new task row (row 14) ---
# Copy formatting (fill color) from E13 (the "Mario" colored cell) onto E14
# so the new "Responsable" cell gets the same visual style used elsewhere
# for Mario.
$ws.Range("E13").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D14").Value = "Validar que se cierra la ventana de selección de promociones despes de elegida una."
$ws.Range("E14").Value = "Mario"
$ws.Range("F14").Value = "Hecho"
$ws.Range("G14").Value = 2
$ws.Range("W14").Value = 2

# --- Update the on-screen selection to reflect where the user ended up ---
$ws.Range("G13").Select()

# --- Re-touch the day-total merged header cells ---
# (Mirrors what Excel itself does internally when the workbook is resaved
# after an edit: the merged-cell registration order for these five ranges
# shifts to the end of the list.)
$mergedHeaderCells = @("AZ4:BA4", "AO4:AP4", "AR4:AS4", "AU4:AV4", "AX4:AY4")
foreach ($addr in $mergedHeaderCells) {
    $ws.Range($addr).UnMerge()
    $ws.Range($addr).Merge()
}

Write-Output "Row 14 updated and new task string added."
